$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "7.18") need the
# columns Text number format applied first, otherwise Excel auto-converts
# the typed text into a numeric value and the literal formatting (trailing
# zeros, etc.) would be lost.
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D9",
    "D11",
    "D14",
    "D18",
    "D19",
    "D21",
    "D23",
    "D24",
    "D26",
    "D27",
    "D29",
    "D30",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D41",
    "D43",
    "D44",
    "D47",
    "D48",
    "D50",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values scraped on
# Mon Jul  8 18:40:55 UTC 2024.
$ws.Range("D2").Value = '56.492.09'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '2.991.09'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '506.72'
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("D6").Value = '137.50'
$ws.Range("E6").Value = '  -0.64%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.429'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '7.18'
$ws.Range("E9").Value = '  -1.75%  '
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("D11").Value = '0.366'
$ws.Range("E11").Value = '  +1.81%  '
$ws.Range("D12").Value = '3.493.61'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("D14").Value = '25.72'
$ws.Range("E14").Value = '  -1.74%  '
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").Value = '56.367.26'
$ws.Range("E16").Value = '  -1.49%  '
$ws.Range("D17").Value = '2.985.65'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '5.99'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").Value = '12.92'
$ws.Range("E19").Value = '  +2.05%  '
$ws.Range("E20").Value = '  +2.19%  '
$ws.Range("D21").Value = '331.95'
$ws.Range("E21").Value = '  +3.53%  '
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").Value = '0.493'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '64.64'
$ws.Range("E24").Value = '  +1.91%  '
$ws.Range("D25").Value = '3.106.29'
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.10%  '
$ws.Range("D27").Value = '0.164'
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("D28").Value = '0.0₃0921'
$ws.Range("E28").Value = '  +2.65%  '
$ws.Range("D29").Value = '6.36'
$ws.Range("E29").Value = '  -3.04%  '
$ws.Range("D30").Value = '6.92'
$ws.Range("E30").Value = '  -2.94%  '
$ws.Range("E31").Value = '  +0.14%  '
$ws.Range("D32").Value = '20.24'
$ws.Range("E32").Value = '  -0.12%  '
$ws.Range("E33").Value = '  -1.26%  '
$ws.Range("D34").Value = '152.95'
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("D35").Value = '4.47'
$ws.Range("E35").Value = '  -2.31%  '
$ws.Range("D36").Value = '5.80'
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("D37").Value = '26.15'
$ws.Range("E37").Value = '  +7.02%  '
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("E39").Value = '  -0.69%  '
$ws.Range("D40").Value = '3.021.10'
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").Value = '36.93'
$ws.Range("E41").Value = '  -2.22%  '
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D43").Value = '3.80'
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("D44").Value = '0.651'
$ws.Range("E44").Value = '  +0.80%  '
$ws.Range("D45").Value = '2.180.57'
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("E46").Value = '  -2.84%  '
$ws.Range("D47").Value = '5.82'
$ws.Range("E47").Value = '  -2.19%  '
$ws.Range("D48").Value = '0.921'
$ws.Range("E48").Value = '  -2.10%  '
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").Value = '19.50'
$ws.Range("E50").Value = '  +1.03%  '
$ws.Range("D51").Value = '0.0850'
$ws.Range("E51").Value = '  -2.89%  '
